$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue "D2" "26.535.88"
Set-TextValue "E2" "  +4.00%  "
Set-TextValue "D3" "1.738.99"
Set-TextValue "E3" "  +4.45%  "
Set-TextValue "D4" "0.9997"
Set-TextValue "E4" "  +0.06%  "
Set-TextValue "D5" "245.11"
Set-TextValue "E5" "  +4.38%  "
Set-TextValue "D6" "1.000"
Set-TextValue "E6" "  -0.04%  "
Set-TextValue "D7" "0.4804"
Set-TextValue "E7" "  +3.88%  "
Set-TextValue "D8" "0.2678"
Set-TextValue "E8" "  +4.12%  "
Set-TextValue "D9" "0.06236"
Set-TextValue "E9" "  +1.93%  "
Set-TextValue "D10" "1.738.37"
Set-TextValue "E10" "  +4.46%  "
Set-TextValue "D11" "0.07133"
Set-TextValue "E11" "  +2.84%  "
Set-TextValue "D12" "15.82"
Set-TextValue "E12" "  +8.33%  "
Set-TextValue "D13" "0.6206"
Set-TextValue "E13" "  +8.33%  "
Set-TextValue "D14" "4.535"
Set-TextValue "E14" "  +4.50%  "
Set-TextValue "D15" "77.05"
Set-TextValue "E15" "  +2.83%  "
Set-TextValue "D16" "0.9997"
Set-TextValue "E16" "  -0.08%  "
Set-TextValue "D17" "26.549.42"
Set-TextValue "E17" "  +4.06%  "
Set-TextValue "D18" "0.9999"
Set-TextValue "E18" "  -0.02%  "
Set-TextValue "D19" "0.000006896"
Set-TextValue "E19" "  +2.72%  "
Set-TextValue "D20" "11.76"
Set-TextValue "E20" "  +3.51%  "
Set-TextValue "D21" "1.962.96"
Set-TextValue "E21" "  +4.50%  "
Set-TextValue "D22" "4.576"
Set-TextValue "E22" "  +4.00%  "
Set-TextValue "D23" "8.921"
Set-TextValue "E23" "  +3.12%  "
Set-TextValue "D24" "5.348"
Set-TextValue "E24" "  +2.44%  "
Set-TextValue "D25" "135.34"
Set-TextValue "E25" "  +0.51%  "
Set-TextValue "E26" "  +3.93%  "
Set-TextValue "D27" "1.814"
Set-TextValue "E27" "  +6.09%  "
Set-TextValue "D28" "1.420"
Set-TextValue "E28" "  +4.21%  "
Set-TextValue "D29" "106.98"
Set-TextValue "E29" "  +3.28%  "
Set-TextValue "D30" "4.000"
Set-TextValue "E30" "  +1.26%  "
Set-TextValue "D31" "3.738"
Set-TextValue "E31" "  +4.07%  "
Set-TextValue "D32" "0.07902"
Set-TextValue "E32" "  +2.74%  "
Set-TextValue "D33" "0.04588"
Set-TextValue "E33" "  +6.44%  "
Set-TextValue "D34" "2.615"
Set-TextValue "E34" "  -0.15%  "
Set-TextValue "D35" "0.9999"
Set-TextValue "E35" "  +6.33%  "
Set-TextValue "D36" "0.6364"
Set-TextValue "E36" "  +6.38%  "
Set-TextValue "D37" "0.9292"
Set-TextValue "E37" "  +1.48%  "
Set-TextValue "D38" "111.99"
Set-TextValue "E38" "  +5.20%  "
Set-TextValue "D39" "1.994"
Set-TextValue "D40" "2.433"
Set-TextValue "E40" "  -1.51%  "
Set-TextValue "D41" "1.005"
Set-TextValue "E41" "  +0.55%  "
Set-TextValue "D42" "0.01515"
Set-TextValue "E42" "  +3.86%  "
Set-TextValue "D43" "5.727"
Set-TextValue "E43" "  +15.31%  "
Set-TextValue "D44" "0.3918"
Set-TextValue "E44" "  +5.82%  "
Set-TextValue "D45" "6.951"
Set-TextValue "E45" "  +13.95%  "
Set-TextValue "D46" "0.1200"
Set-TextValue "E46" "  +8.18%  "
Set-TextValue "D47" "0.05332"
Set-TextValue "E47" "  +1.49%  "
Set-TextValue "B48" "Elrond"
Set-TextValue "C48" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue "D48" "30.85"
Set-TextValue "E48" "  +3.17%  "
Set-TextValue "B49" "EnergySwap"
Set-TextValue "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "7.885"
Set-TextValue "E49" "  +4.46%  "
Set-TextValue "D50" "1.256"
Set-TextValue "E50" "  +6.00%  "
Set-TextValue "D51" "0.3448"
Set-TextValue "E51" "  +4.57%  "
